# Insert a new row above row 84 (shifts rows 84:148 down to 85:149,
# carrying formatting/styles with them, same as a native Excel row insert).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(84).Insert()

# Populate the newly inserted row 84 with the new weekly record.
# Columns A,B,C,E,F,G,H,I,J,K,L,Q,R,T mirror the row immediately below it
# (the old row 84, now shifted to row 85) since those fields are unchanged;
# only D (Fecha), M (Volumen), N (Precio minimo), O (Precio maximo),
# P (Precio promedio ponderado) and S (Precio $/Kg) carry new values.
$ws.Range("A84").Value = 7
$ws.Range("B84").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C84").Value = "Ñuble"
$ws.Range("D84").Value = 45096
$ws.Range("E84").Value = 16
$ws.Range("F84").Value = "Fruta"
$ws.Range("G84").Value = 100108
$ws.Range("H84").Value = "Tropicales y subtropicales"
$ws.Range("I84").Value = 100108002
$ws.Range("J84").Value = "Mango"
$ws.Range("K84").Value = "Sin especificar"
$ws.Range("L84").Value = "Primera"
$ws.Range("M84").Value = 50
$ws.Range("N84").Value = 9000
$ws.Range("O84").Value = 9000
$ws.Range("P84").Value = 9000
$ws.Range("Q84").Value = "$/bandeja 4 kilos"
$ws.Range("R84").Value = "Perú"
$ws.Range("S84").Value = 2250
$ws.Range("T84").Value = 4
